# "new script for 2015"
# The raw-data sheet ("Sheet1") had a block of cells holding the literal
# string "<Null>" (shared-string placeholder) in columns where no bumble
# bee of a given species was observed within that buffer radius. The new
# script that generates this workbook instead emits a numeric 0 for those
# same cells, so replicate that data correction here.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ranges = @(
    "I4:N4",
    "I9:L9",
    "I10:L10",
    "I11:N11",
    "I13:J13",
    "I14:R14",
    "I15:T15",
    "I16:J16",
    "I17:P17",
    "I18:R18",
    "I19:P19",
    "I20:L20"
)

foreach ($rng in $ranges) {
    $ws.Range($rng).Value = 0
}

# Restore the view state: the workbook now opens on Sheet1 (rather than
# Sheet9), with the selection left on K20 ...
$ws.Activate()
$ws.Range("K20").Select()

# ... and the "5km" sheet's lingering multi-cell selection collapses back
# to a single cell.
$ws5km = $wb.Worksheets.Item("5km")
$ws5km.Range("AL2").Select()

# Sheet1 is the one that should end up tabSelected/active, so re-activate
# it last.
$ws.Activate()
